# Update attendance workbook: mark "Real" attendance (D/E) and "Absent" (H)
# for each date row, and flag the first date as "Invalid" (G) too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid + Absent
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Rows 4-6: Real attendance (Total + Real columns)
$ws.Range("D4:E4").Value = 1
$ws.Range("D5:E5").Value = 1
$ws.Range("D6:E6").Value = 1

# Rows 7-9: Absent
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1

# Rows 10-11: Real attendance
$ws.Range("D10:E10").Value = 1
$ws.Range("D11:E11").Value = 1

# Rows 12-18: Absent
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
